# Adds the 2024/11/08 data column (BI) to the "合成確率" sheet, mirroring the
# existing per-day columns (B..BH): a text date header in row 1 plus 52 numeric
# values in rows 2-53, each carrying one of the three existing cell styles
# (s=1 plain, s=2 yellow fill, s=3 light-blue fill).
#
# -4122 == xlPasteFormats (Excel's PasteSpecial "Paste: Formats" option). We
# reuse it to copy the exact cell style (font + fill) from a representative
# cell that already has the desired style index, instead of rebuilding the
# style from scratch (which would otherwise create new, slightly different
# style/font entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BI (61) width, matching neighboring data columns (XML width=12) ---
$ws.Columns.Item(61).ColumnWidth = 11.17

# --- Header cell BI1: new date label "2024/11/08" as text, styled like BH1 ---
$ws.Range("BI1").Value = "'2024/11/08"
$ws.Range("BH1").Copy()
$ws.Range("BI1").PasteSpecial(-4122)

# --- Data cells BI2:BI53 (numeric values) ---
$ws.Range("BI2").Value = 153.8
$ws.Range("BI3").Value = 179.8
$ws.Range("BI4").Value = 216.1
$ws.Range("BI5").Value = 165.8
$ws.Range("BI6").Value = 331
$ws.Range("BI7").Value = 128.6
$ws.Range("BI8").Value = 156.6
$ws.Range("BI9").Value = 155.6
$ws.Range("BI10").Value = 143.1
$ws.Range("BI11").Value = 166.8
$ws.Range("BI12").Value = 147.3
$ws.Range("BI13").Value = 454
$ws.Range("BI14").Value = 176.4
$ws.Range("BI15").Value = 129.5
$ws.Range("BI16").Value = 154.2
$ws.Range("BI17").Value = 122.3
$ws.Range("BI18").Value = 239
$ws.Range("BI19").Value = 154.6
$ws.Range("BI20").Value = 138.3
$ws.Range("BI21").Value = 174.4
$ws.Range("BI22").Value = 156.1
$ws.Range("BI23").Value = 150
$ws.Range("BI24").Value = 178.1
$ws.Range("BI25").Value = 155.1
$ws.Range("BI26").Value = 152.1
$ws.Range("BI27").Value = 193.5
$ws.Range("BI28").Value = 149
$ws.Range("BI29").Value = 183.1
$ws.Range("BI30").Value = 181.1
$ws.Range("BI31").Value = 161.5
$ws.Range("BI32").Value = 108.6
$ws.Range("BI33").Value = 156.3
$ws.Range("BI34").Value = 105.3
$ws.Range("BI35").Value = 178.7
$ws.Range("BI36").Value = 129.5
$ws.Range("BI37").Value = 117.8
$ws.Range("BI38").Value = 138.6
$ws.Range("BI39").Value = 157
$ws.Range("BI40").Value = 219.2
$ws.Range("BI41").Value = 141.1
$ws.Range("BI42").Value = 182.4
$ws.Range("BI43").Value = 129.1
$ws.Range("BI44").Value = 151.5
$ws.Range("BI45").Value = 129.1
$ws.Range("BI46").Value = 176.6
$ws.Range("BI47").Value = 151.2
$ws.Range("BI48").Value = 211.1
$ws.Range("BI49").Value = 225
$ws.Range("BI50").Value = 147.3
$ws.Range("BI51").Value = 139.8
$ws.Range("BI52").Value = 154.2
$ws.Range("BI53").Value = 165.7

# --- Apply matching cell styles (font/fill) per diff, via copy/paste-format from representative cells ---
$ws.Range("A2").Copy()
$ws.Range("BI2").PasteSpecial(-4122)
$ws.Range("BI3").PasteSpecial(-4122)
$ws.Range("BI4").PasteSpecial(-4122)
$ws.Range("BI5").PasteSpecial(-4122)
$ws.Range("BI6").PasteSpecial(-4122)
$ws.Range("BI8").PasteSpecial(-4122)
$ws.Range("BI9").PasteSpecial(-4122)
$ws.Range("BI10").PasteSpecial(-4122)
$ws.Range("BI11").PasteSpecial(-4122)
$ws.Range("BI12").PasteSpecial(-4122)
$ws.Range("BI13").PasteSpecial(-4122)
$ws.Range("BI14").PasteSpecial(-4122)
$ws.Range("BI16").PasteSpecial(-4122)
$ws.Range("BI18").PasteSpecial(-4122)
$ws.Range("BI19").PasteSpecial(-4122)
$ws.Range("BI21").PasteSpecial(-4122)
$ws.Range("BI22").PasteSpecial(-4122)
$ws.Range("BI23").PasteSpecial(-4122)
$ws.Range("BI24").PasteSpecial(-4122)
$ws.Range("BI25").PasteSpecial(-4122)
$ws.Range("BI26").PasteSpecial(-4122)
$ws.Range("BI27").PasteSpecial(-4122)
$ws.Range("BI28").PasteSpecial(-4122)
$ws.Range("BI29").PasteSpecial(-4122)
$ws.Range("BI30").PasteSpecial(-4122)
$ws.Range("BI31").PasteSpecial(-4122)
$ws.Range("BI33").PasteSpecial(-4122)
$ws.Range("BI35").PasteSpecial(-4122)
$ws.Range("BI39").PasteSpecial(-4122)
$ws.Range("BI40").PasteSpecial(-4122)
$ws.Range("BI41").PasteSpecial(-4122)
$ws.Range("BI42").PasteSpecial(-4122)
$ws.Range("BI44").PasteSpecial(-4122)
$ws.Range("BI46").PasteSpecial(-4122)
$ws.Range("BI47").PasteSpecial(-4122)
$ws.Range("BI48").PasteSpecial(-4122)
$ws.Range("BI49").PasteSpecial(-4122)
$ws.Range("BI50").PasteSpecial(-4122)
$ws.Range("BI52").PasteSpecial(-4122)
$ws.Range("BI53").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("BI17").PasteSpecial(-4122)
$ws.Range("BI32").PasteSpecial(-4122)
$ws.Range("BI34").PasteSpecial(-4122)
$ws.Range("BI37").PasteSpecial(-4122)
$ws.Range("N2").Copy()
$ws.Range("BI7").PasteSpecial(-4122)
$ws.Range("BI15").PasteSpecial(-4122)
$ws.Range("BI20").PasteSpecial(-4122)
$ws.Range("BI36").PasteSpecial(-4122)
$ws.Range("BI38").PasteSpecial(-4122)
$ws.Range("BI43").PasteSpecial(-4122)
$ws.Range("BI45").PasteSpecial(-4122)
$ws.Range("BI51").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Output "applied BI column for 2024/11/08"
